$d = $word.ActiveDocument

# --- Paragraph 1: bold heading "Web Server Literature Review" ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Text = "Web Server Literature Review"
$p1b = $d.Paragraphs.Item(1)
$r1b = $p1b.Range
$r1b.Font.Bold = $true

# --- Paragraph 2: quote paragraph, paragraph-mark colored red, run left plain ---
$quote = "The article explains that " + [char]0x201C + "We have used our implementations to carry out a bottleneck characterization of the benchmarks. Different benchmarks show different bottlenecks: the database CPU for the online bookstore, and the Web server CPU for the auction site and the bulletin board. Complex queries cause the database CPU to the bottleneck for the online bookstore. In contrast, the queries for the other applications are simpler" + [char]0x201D + " [13]."

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Text = $quote

# Color the whole paragraph (text + paragraph mark) red -- this paints both
# the run and the paragraph-mark rPr.
$p2b = $d.Paragraphs.Item(2)
$r2full = $p2b.Range
$r2full.Font.Color = 255

# Re-create the run text (delete + re-insert) so the run itself goes back to
# being unformatted while the paragraph mark keeps the red color.
$p2c = $d.Paragraphs.Item(2)
$r2text = $p2c.Range
[void]$r2text.MoveEnd(1, -1)
$r2text.Delete()
$p2d = $d.Paragraphs.Item(2)
$r2d = $p2d.Range
$r2d.InsertBefore($quote)

# --- Remove the two now-redundant empty Style3 paragraphs (old #3 and #4) ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Delete()
$p3again = $d.Paragraphs.Item(3)
$p3again.Range.Delete()

Write-Output "done"
